$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'282.02"
$ws.Range("D3").Value = "'20.97"
$ws.Range("D4").Value = "'6.245"
$ws.Range("D5").Value = "'0.06160"
$ws.Range("D7").Value = "'6.556"
$ws.Range("D8").Value = "'1.477"
$ws.Range("D9").Value = "'0.8156"
$ws.Range("D10").Value = "'0.01386"
$ws.Range("D11").Value = "'0.1638"
$ws.Range("D12").Value = "'0.08321"
$ws.Range("D13").Value = "'0.03538"
$ws.Range("D14").Value = "'0.03149"
$ws.Range("D15").Value = "'0.09140"
$ws.Range("D16").Value = "'3.727"
$ws.Range("D17").Value = "'0.001639"
$ws.Range("D18").Value = "'0.04662"
$ws.Range("D19").Value = "'0.006441"
$ws.Range("D20").Value = "'0.006191"
$ws.Range("D21").Value = "'0.001067"
$ws.Range("D23").Value = "'3.815"
$ws.Range("D25").Value = "'0.3374"
$ws.Range("D26").Value = "'0.1249"
$ws.Range("D40").Value = "'0.04674"
$ws.Range("D41").Value = "'0.007121"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1106"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003518"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.01147"
$ws.Range("D45").Value = "'0.00006259"
$ws.Range("D47").Value = "'0.9992"
$ws.Range("D48").Value = "'0.002916"
